# Abstract rewrite: replace the placeholder "Abstract" paragraph with the
# USFWS / Battle Creek abstract copy, switch the paragraph to the Aptos
# font (carried on both the paragraph mark and every run), and leave two
# blank paragraphs trailing the abstract text.
#
# The abstract text is authored as a sequence of separate runs (this is
# how it came out of the source edit - e.g. pasted/typed in several
# passes), so rather than fight Word's "identical formatting -> merge
# into one run" behavior with repeated Selection.TypeText/Find.Execute
# calls, we build the exact run-by-run Open XML for the paragraph and
# drop it in with Range.InsertXML, which inserts the markup verbatim.

$d = $word.ActiveDocument

# Each entry is the run's literal text plus whether it needs
# xml:space="preserve" (i.e. has leading/trailing whitespace).
$runs = @(
    @{ Text = "The "; Preserve = $true },
    @{ Text = "U.S. Fish and Wildlife Service"; Preserve = $false },
    @{ Text = " ("; Preserve = $true },
    @{ Text = "USFWS"; Preserve = $false },
    @{ Text = ") "; Preserve = $true },
    @{ Text = "collects data on adult salmonids on "; Preserve = $true },
    @{ Text = "Battle Creek"; Preserve = $false },
    @{ Text = ". Data is collected annually via "; Preserve = $true },
    @{ Text = "snorkel"; Preserve = $false },
    @{ Text = " surveys; video camera systems collect data on upstream passage 24 hours a day, 7 days a week from "; Preserve = $true },
    @{ Text = "March to November"; Preserve = $false },
    @{ Text = ". Data from this monitoring is used to estimate adult escapement (upstream passage) abundance and timing, spawner abundance, and other important metrics for adult salmonids in the watershed. These"; Preserve = $false },
    @{ Text = " data will also be used to inform the development of a juvenile production estimate (JPE) for spring-run Chinook salmon in the Sacramento River Watershed"; Preserve = $true },
    @{ Text = "."; Preserve = $false }
)

$rPr = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/></w:rPr>'

$runsXml = ""
foreach ($run in $runs) {
    if ($run.Preserve) { $spaceAttr = ' xml:space="preserve"' } else { $spaceAttr = '' }
    $runsXml += "<w:r>$rPr<w:t$spaceAttr>$($run.Text)</w:t></w:r>"
}

# Paragraph 1: the abstract text, Aptos on the paragraph mark (pPr/rPr) too.
$para1 = "<w:p><w:pPr>$rPr</w:pPr>$runsXml</w:p>"
# Paragraph 2: blank, but still carries the Aptos paragraph-mark font.
$para2 = "<w:p><w:pPr>$rPr</w:pPr></w:p>"
# Paragraph 3: a plain trailing blank paragraph.
$para3 = "<w:p/>"
# Extra trailing <w:p/> is absorbed into the end-of-range paragraph mark
# rather than becoming a fourth visible paragraph - it's required so
# paragraph 3 above actually lands as its own paragraph.
$trailingMarker = "<w:p/>"

$bodyXml = "<w:body>$para1$para2$para3$trailingMarker</w:body>"

$packageXml = '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    $bodyXml +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

# Replace the whole document body (the single "Abstract" paragraph) with
# the rebuilt abstract content in one shot.
$d.Content.InsertXML($packageXml)
